$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1971
    $ws.Range("I2").Value = "https://show.bilibili.com/platform/detail.html?id=80053"

    $ws.Range("F3").Value = 269
    $ws.Range("I3").Value = "https://show.bilibili.com/platform/detail.html?id=79124"

    $ws.Range("F4").Value = 268
    $ws.Range("I4").Value = "https://show.bilibili.com/platform/detail.html?id=79002"

    $ws.Range("F5").Value = 9052
    $ws.Range("I5").Value = "https://show.bilibili.com/platform/detail.html?id=80084"

    $ws.Range("F6").Value = 10442
    $ws.Range("I6").Value = "https://show.bilibili.com/platform/detail.html?id=80426"

    $ws.Range("I7").Value = "https://show.bilibili.com/platform/detail.html?id=80398"

    $ws.Range("I8").Value = "https://show.bilibili.com/platform/detail.html?id=80528"

    $ws.Range("F9").Value = 661
    $ws.Range("I9").Value = "https://show.bilibili.com/platform/detail.html?id=80504"

    $ws.Range("F10").Value = 89
    $ws.Range("I10").Value = "https://show.bilibili.com/platform/detail.html?id=80248"

    $ws.Range("F11").Value = 9394
    $ws.Range("I11").Value = "https://show.bilibili.com/platform/detail.html?id=79303"

    $ws.Range("F12").Value = 12
    $ws.Range("I12").Value = "https://show.bilibili.com/platform/detail.html?id=81044"

    $ws.Range("F13").Value = 2411
    $ws.Range("I13").Value = "https://show.bilibili.com/platform/detail.html?id=79333"

    $ws.Range("F14").Value = 23
    $ws.Range("I14").Value = "https://show.bilibili.com/platform/detail.html?id=80635"

    $ws.Range("F15").Value = 55
    $ws.Range("I15").Value = "https://show.bilibili.com/platform/detail.html?id=80789"

    $ws.Range("F16").Value = 351
    $ws.Range("I16").Value = "https://show.bilibili.com/platform/detail.html?id=78666"

    $ws.Range("F17").Value = 10667
    $ws.Range("I17").Value = "https://show.bilibili.com/platform/detail.html?id=79789"

    $ws.Range("F18").Value = 10598
    $ws.Range("I18").Value = "https://show.bilibili.com/platform/detail.html?id=77196"

    $ws.Range("F19").Value = 2
    $ws.Range("I19").Value = "https://show.bilibili.com/platform/detail.html?id=81118"

    $ws.Range("F20").Value = 3
    $ws.Range("I20").Value = "https://show.bilibili.com/platform/detail.html?id=81116"

    $ws.Range("F21").Value = 2
    $ws.Range("I21").Value = "https://show.bilibili.com/platform/detail.html?id=81119"

    $ws.Range("I22").Value = "https://show.bilibili.com/platform/detail.html?id=81100"

    $ws.Range("F23").Value = 7
    $ws.Range("I23").Value = "https://show.bilibili.com/platform/detail.html?id=81120"

    $ws.Range("I24").Value = "https://show.bilibili.com/platform/detail.html?id=81114"

}